$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2024 April (row 7)
$ws.Range("B7").Value = 55
$ws.Range("C7").Value = 12

# 2024 August (row 11)
$ws.Range("B11").Value = 37
$ws.Range("C11").Value = 19

# 2024 September (row 12)
$ws.Range("B12").Value = 28
$ws.Range("C12").Value = 20

# 2024 October (row 13)
$ws.Range("B13").Value = 9
$ws.Range("C13").Value = 16
